# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with the
# freshly scraped rows for Linea 141 (commit: "Horarios actualizados Linea 141 - 631").
# Each sheet keeps its header block (title / last-updated / total-filas) in rows 1-3,
# a blank row 4, column headers in row 5, and data rows sorted by Hora_Llegada from
# row 6 onward. New rows were scraped, shifting/overwriting the data block and
# growing the sheet's used range.

$wb = $excel.ActiveWorkbook

$rowsSheet1 = @(
    @("03:45:25", "03:47", "14_ABASTO", 2, "LP1912"),
    @("03:45:25", "04:01", "81_EL PELIGRO", 16, "LP1912"),
    @("03:45:25", "04:46", "215A_EL PATO", 61, "LP1912"),
    @("03:45:25", "04:53", "11_ETCHEVERRY", 68, "LP1912"),
    @("04:56:49", "05:13", "14_ABASTO", 17, "LP1912"),
    @("03:45:25", "05:16", "17_ROMERO", 91, "LP1912"),
    @("04:45:05", "05:16", "14_ABASTO", 31, "LP1912"),
    @("03:45:25", "05:22", "23_HERNANDEZ", 97, "LP1912"),
    @("05:26:08", "05:28", "14_ABASTO", 2, "LP1912"),
    @("04:18:02", "05:34", "14_ABASTO", 76, "LP1912"),
    @("03:45:25", "05:34", "215B_EL PATO", 109, "LP1912"),
    @("04:18:02", "05:35", "215B_EL PATO", 77, "LP1912"),
    @("03:45:25", "05:37", "14_ABASTO", 112, "LP1912"),
    @("04:18:02", "05:46", "15_ABASTO", 88, "LP1912"),
    @("04:45:05", "06:04", "16_SANTA ANA", 79, "LP1912"),
    @("04:18:02", "06:05", "16_SANTA ANA", 107, "LP1912"),
    @("04:56:49", "06:11", "215A_EL PATO", 75, "LP1912"),
    @("04:18:02", "06:12", "215A_EL PATO", 114, "LP1912"),
    @("04:18:02", "06:14", "225_HARAS DEL SUR", 116, "LP1912"),
    @("04:45:05", "06:21", "26_HERNANDEZ", 96, "LP1912"),
    @("06:25:43", "06:26", "86_EST CHICA-ESC AGRARIA", 1, "LP1912"),
    @("04:45:05", "06:27", "23_HERNANDEZ", 102, "LP1912"),
    @("06:25:43", "06:28", "23_HERNANDEZ", 3, "LP1912"),
    @("04:56:49", "06:29", "86_EST CHICA-ESC AGRARIA", 93, "LP1912"),
    @("04:45:05", "06:30", "86_EST CHICA-ESC AGRARIA", 105, "LP1912"),
    @("04:45:05", "06:31", "16_SANTA ANA", 106, "LP1912"),
    @("05:55:25", "06:44", "26_HERNANDEZ", 49, "LP1912"),
    @("04:45:05", "06:44", "225_C ROCA-H SUR", 119, "LP1912"),
    @("04:56:49", "06:46", "215C_EL PATO", 110, "LP1912"),
    @("05:26:08", "06:47", "215C_EL PATO", 81, "LP1912"),
    @("05:55:25", "06:59", "14_ABASTO", 64, "LP1912"),
    @("05:26:08", "07:00", "14_ABASTO", 94, "LP1912"),
    @("06:25:43", "07:01", "16_SANTA ANA", 36, "LP1912"),
    @("05:55:25", "07:04", "23_HERNANDEZ", 69, "LP1912"),
    @("05:26:08", "07:05", "23_HERNANDEZ", 99, "LP1912"),
    @("05:26:08", "07:05", "15_ABASTO", 99, "LP1912"),
    @("05:26:08", "07:06", "10_OLMOS", 100, "LP1912"),
    @("05:26:08", "07:07", "225_GOMEZ", 101, "LP1912"),
    @("05:26:08", "07:11", "215A_EL PATO", 105, "LP1912"),
    @("06:55:02", "07:12", "215A_EL PATO", 17, "LP1912"),
    @("06:25:43", "07:14", "26_HERNANDEZ", 49, "LP1912"),
    @("05:55:25", "07:15", "11_ETCHEVERRY", 80, "LP1912"),
    @("05:26:08", "07:16", "11_ETCHEVERRY", 110, "LP1912"),
    @("06:55:02", "07:17", "16_SANTA ANA", 22, "LP1912"),
    @("05:26:08", "07:21", "26_HERNANDEZ", 115, "LP1912"),
    @("05:26:08", "07:23", "10_OLMOS", 117, "LP1912"),
    @("05:55:25", "07:30", "10_OLMOS", 95, "LP1912"),
    @("05:55:25", "07:31", "16_SANTA ANA", 96, "LP1912"),
    @("05:55:25", "07:31", "11_ETCHEVERRY", 96, "LP1912"),
    @("06:55:02", "07:32", "16_SANTA ANA", 37, "LP1912"),
    @("05:55:25", "07:32", "84_COLONIA URQUIZA-ESC 49", 97, "LP1912"),
    @("06:55:02", "07:32", "11_ETCHEVERRY", 37, "LP1912"),
    @("05:55:25", "07:36", "27_EL RETIRO", 101, "LP1912"),
    @("06:55:02", "07:37", "27_EL RETIRO", 42, "LP1912"),
    @("05:55:25", "07:39", "10_OLMOS", 104, "LP1912"),
    @("05:55:25", "07:47", "14_ABASTO", 112, "LP1912"),
    @("06:55:02", "07:48", "14_ABASTO", 53, "LP1912"),
    @("05:55:25", "07:51", "215D_EL PATO", 116, "LP1912"),
    @("06:55:02", "07:52", "215D_EL PATO", 57, "LP1912"),
    @("06:25:43", "08:01", "23_HERNANDEZ", 96, "LP1912"),
    @("06:55:02", "08:03", "23_HERNANDEZ", 68, "LP1912"),
    @("06:25:43", "08:12", "15_ABASTO", 107, "LP1912"),
    @("06:55:02", "08:21", "26_HERNANDEZ", 86, "LP1912"),
    @("06:25:43", "08:22", "16_P MOR-SANTA ANA", 117, "LP1912"),
    @("06:55:02", "08:23", "16_P MOR-SANTA ANA", 88, "LP1912"),
    @("06:25:43", "08:23", "215B_EL PATO", 118, "LP1912"),
    @("06:55:02", "08:27", "84_COLONIA URQUIZA-ESC 49", 92, "LP1912"),
    @("06:55:02", "08:42", "81_EL PELIGRO", 107, "LP1912"),
    @("06:55:02", "08:54", "17_ROMERO", 119, "LP1912"),
)

$rowsSheet2 = @(
    @("03:45:25", "04:46", "215A_EL PATO", 61, "LP1912"),
    @("03:45:25", "05:34", "215B_EL PATO", 109, "LP1912"),
    @("04:18:02", "05:35", "215B_EL PATO", 77, "LP1912"),
    @("04:56:49", "06:11", "215A_EL PATO", 75, "LP1912"),
    @("04:18:02", "06:12", "215A_EL PATO", 114, "LP1912"),
    @("04:56:49", "06:46", "215C_EL PATO", 110, "LP1912"),
    @("05:26:08", "06:47", "215C_EL PATO", 81, "LP1912"),
    @("05:26:08", "07:11", "215A_EL PATO", 105, "LP1912"),
    @("06:55:02", "07:12", "215A_EL PATO", 17, "LP1912"),
    @("05:55:25", "07:51", "215D_EL PATO", 116, "LP1912"),
    @("06:55:02", "07:52", "215D_EL PATO", 57, "LP1912"),
    @("06:25:43", "08:23", "215B_EL PATO", 118, "LP1912"),
)

$rowsSheet3 = @(
    @("04:56:49", "05:43", "215A_LA PLATA", 47, "L6173"),
    @("03:45:25", "05:44", "215A_LA PLATA", 119, "L6173"),
    @("04:56:49", "06:08", "215A_LA PLATA", 72, "L6173"),
    @("04:18:02", "06:09", "215A_LA PLATA", 111, "L6173"),
    @("04:56:49", "06:32", "215C_LA PLATA", 96, "L6203"),
    @("04:45:05", "06:33", "215C_LA PLATA", 108, "L6203"),
    @("06:25:43", "06:59", "215B_LP-P MOR-1 Y 57", 34, "L6173"),
    @("05:26:08", "07:00", "215B_LP-P MOR-1 Y 57", 94, "L6173"),
    @("05:55:25", "07:35", "215A_LA PLATA", 100, "L6173"),
    @("06:25:43", "07:39", "215A_LA PLATA", 74, "L6173"),
    @("06:55:02", "07:42", "215A_LA PLATA", 47, "L6173"),
    @("06:25:43", "08:06", "215C_LA PLATA", 101, "L6203"),
    @("06:55:02", "08:07", "215C_LA PLATA", 72, "L6203"),
    @("06:55:02", "08:36", "215A_LA PLATA", 101, "L6173"),
)


function Update-ScheduleSheet($ws, $lastUpdated, $totalFilas, $rows) {
    $ws.Cells.Item(2,1).Value = $lastUpdated
    $ws.Cells.Item(3,1).Value = $totalFilas

    $startRow = 6
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $rows[$i]
        $rowNum = $startRow + $i
        $ws.Cells.Item($rowNum,1).Value = $r[0]
        $ws.Cells.Item($rowNum,2).Value = $r[1]
        $ws.Cells.Item($rowNum,3).Value = $r[2]
        $ws.Cells.Item($rowNum,4).Value = $r[3]
        $ws.Cells.Item($rowNum,5).Value = $r[4]
    }
}

$lastUpdated = "Última actualización: 06:55:02"

$ws1 = $wb.Worksheets.Item("LP1912")
Update-ScheduleSheet $ws1 $lastUpdated "Total filas: 69" $rowsSheet1

$ws2 = $wb.Worksheets.Item("LP1912-215")
Update-ScheduleSheet $ws2 $lastUpdated "Total filas: 12" $rowsSheet2

$ws3 = $wb.Worksheets.Item("6203-6173")
Update-ScheduleSheet $ws3 $lastUpdated "Total filas: 14" $rowsSheet3

"Updated {0} / {1} / {2} rows" -f $rowsSheet1.Count, $rowsSheet2.Count, $rowsSheet3.Count
